# Region III_HEALTH.xlsx - update to most recent status/accomplishment report
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1. Header row updates
# ---------------------------------------------------------------------------
# I1 / L1 text now in upper case
$ws.Range("I1").Value = "TOTAL PHYSICAL TARGET"
$ws.Range("L1").Value = "BATCH"

# Move the existing "Status as of ..." header from AA1 out to AF1 to make room
# for five new "No. of Sites ..." summary columns.
$ws.Range("AF1").Value = $ws.Range("AA1").Value2

# Give the five new header cells (AA1:AE1) the same look as the rest of the
# header row (bold / centered / bordered) by copying the formatting of an
# existing header cell, then fill in their text.
$ws.Range("A1").Copy($ws.Range("AA1:AE1"))
$ws.Range("AA1").Value = "No. of Sites Reverted"
$ws.Range("AB1").Value = "No. of Sites Not yet started"
$ws.Range("AC1").Value = "No. of Sites Under Procurement"
$ws.Range("AD1").Value = "No. of Sites On Going"
$ws.Range("AE1").Value = "No. of Sites Completed"

# ---------------------------------------------------------------------------
# 2. Clear placeholder "-" values in columns I and L for rows 2-17
# ---------------------------------------------------------------------------
$ws.Range("I2:I17").ClearContents()
$ws.Range("L2:L17").ClearContents()

# ---------------------------------------------------------------------------
# 3. Move the dropdown data validation from AA2:AA36 to AF2:AF36
# ---------------------------------------------------------------------------
$ws.Range("AA2:AA36").Validation.Delete()
$ws.Range("AF2:AF36").Validation.Add(3, 1, 1, "=DropdownOptions!`$A`$1:`$A`$7")
